$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.955205556636713
$ws.Cells.Item(2, 3).Value = 0.1864482854026051
$ws.Cells.Item(2, 4).Value = 0.186891463224228
$ws.Cells.Item(2, 5).Value = 0.1508914283903522
$ws.Cells.Item(2, 6).Value = 1.292313883070349
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 10).Value = 0.1618723921172318
$ws.Cells.Item(2, 13).Value = 0.366599744145077
$ws.Cells.Item(2, 15).Value = 3.042813547782231

# Row 3
$ws.Cells.Item(3, 2).Value = 0.8538813269291836
$ws.Cells.Item(3, 3).Value = 0.1626727339523768
$ws.Cells.Item(3, 4).Value = 0.1846395855147591
$ws.Cells.Item(3, 5).Value = 0.150854307453173
$ws.Cells.Item(3, 6).Value = 1.298760081499694
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 10).Value = 0.1630871668864025
$ws.Cells.Item(3, 13).Value = 0.3404761391511713
$ws.Cells.Item(3, 15).Value = 3.066688282738795

# Row 4
$ws.Cells.Item(4, 2).Value = 0.7916388058351345
$ws.Cells.Item(4, 3).Value = 0.1480177641776663
$ws.Cells.Item(4, 4).Value = 0.1833096438735495
$ws.Cells.Item(4, 5).Value = 0.1508972905513488
$ws.Cells.Item(4, 6).Value = 1.303623361375685
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 10).Value = 0.1639225429714628
$ws.Cells.Item(4, 13).Value = 0.3244995179948233
$ws.Cells.Item(4, 15).Value = 3.083719490105977

# Row 5
$ws.Cells.Item(5, 2).Value = 0.7662685286573208
$ws.Cells.Item(5, 3).Value = 0.1420318546391002
$ws.Cells.Item(5, 4).Value = 0.1827810063810773
$ws.Cells.Item(5, 5).Value = 0.1509313656494591
$ws.Cells.Item(5, 6).Value = 1.305832574178844
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 10).Value = 0.1642854613743268
$ws.Cells.Item(5, 13).Value = 0.3180052103878168
$ws.Cells.Item(5, 15).Value = 3.09125526657192

# Row 6
$ws.Cells.Item(6, 2).Value = 0.7620554946846596
$ws.Cells.Item(6, 3).Value = 0.1410370707646962
$ws.Cells.Item(6, 4).Value = 0.1826940332928331
$ws.Cells.Item(6, 5).Value = 0.1509380244964333
$ws.Cells.Item(6, 6).Value = 1.306213138403521
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 10).Value = 0.164347082116894
$ws.Cells.Item(6, 13).Value = 0.316927830377594
$ws.Cells.Item(6, 15).Value = 3.092542507885128

# Row 7
$ws.Cells.Item(7, 2).Value = 0.7912966756640003
$ws.Cells.Item(7, 3).Value = 0.1479370918652023
$ws.Cells.Item(7, 4).Value = 0.1833024604425901
$ws.Cells.Item(7, 5).Value = 0.1508976830268907
$ws.Cells.Item(7, 6).Value = 1.303652235253765
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 10).Value = 0.1639273463438542
$ws.Cells.Item(7, 13).Value = 0.3244118670342999
$ws.Cells.Item(7, 15).Value = 3.083818710779695

# Row 8
$ws.Cells.Item(8, 2).Value = 0.9202757663928764
$ws.Cells.Item(8, 3).Value = 0.1782624499438725
$ws.Cells.Item(8, 4).Value = 0.1861041130960004
$ws.Cells.Item(8, 5).Value = 0.1508649886112252
$ws.Cells.Item(8, 6).Value = 1.294348484180915
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 10).Value = 0.1622726650361912
$ws.Cells.Item(8, 13).Value = 0.3575793765504471
$ws.Cells.Item(8, 15).Value = 3.050552612770502

# Row 9
$ws.Cells.Item(9, 2).Value = 1.172926426550589
$ws.Cells.Item(9, 3).Value = 0.2372676085363707
$ws.Cells.Item(9, 4).Value = 0.1920139044422484
$ws.Cells.Item(9, 5).Value = 0.1513220894922398
$ws.Cells.Item(9, 6).Value = 1.283299597575862
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 10).Value = 0.159738469880061
$ws.Cells.Item(9, 13).Value = 0.4231114079752629
$ws.Cells.Item(9, 15).Value = 3.004187772235298

# Row 10
$ws.Cells.Item(10, 2).Value = 1.358336643231439
$ws.Cells.Item(10, 3).Value = 0.2803234656453526
$ws.Cells.Item(10, 4).Value = 0.1966064647408672
$ws.Cells.Item(10, 5).Value = 0.1519749930517094
$ws.Cells.Item(10, 6).Value = 1.279587958492499
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 10).Value = 0.1583105886951266
$ws.Cells.Item(10, 13).Value = 0.4715456480627509
$ws.Cells.Item(10, 15).Value = 2.981697548245563

# Row 11
$ws.Cells.Item(11, 2).Value = 1.442630677157752
$ws.Cells.Item(11, 3).Value = 0.2998441364315454
$ws.Cells.Item(11, 4).Value = 0.1987496094513403
$ws.Cells.Item(11, 5).Value = 0.1523407402143917
$ws.Cells.Item(11, 6).Value = 1.278860405847212
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 10).Value = 0.1577554211149277
$ws.Cells.Item(11, 13).Value = 0.4936401315067229
$ws.Cells.Item(11, 15).Value = 2.973994669067594

# Row 12
$ws.Cells.Item(12, 2).Value = 1.474542391286775
$ws.Cells.Item(12, 3).Value = 0.3072263721833508
$ws.Cells.Item(12, 4).Value = 0.1995688646867109
$ws.Cells.Item(12, 5).Value = 0.152489108876793
$ws.Cells.Item(12, 6).Value = 1.278723391460758
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 10).Value = 0.1575587794229989
$ws.Cells.Item(12, 13).Value = 0.5020153049565863
$ws.Cells.Item(12, 15).Value = 2.971442519217248

# Row 13
$ws.Cells.Item(13, 2).Value = 1.467670037867094
$ws.Cells.Item(13, 3).Value = 0.3056369182773722
$ws.Cells.Item(13, 4).Value = 0.1993920823405517
$ws.Cells.Item(13, 5).Value = 0.1524567164007671
$ws.Cells.Item(13, 6).Value = 1.278746735332064
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 10).Value = 0.1576005251333328
$ws.Cells.Item(13, 13).Value = 0.5002111905917417
$ws.Cells.Item(13, 15).Value = 2.971975927777379

# Row 14
$ws.Cells.Item(14, 2).Value = 1.445256253642469
$ws.Cells.Item(14, 3).Value = 0.3004516763347169
$ws.Cells.Item(14, 4).Value = 0.1988168562804589
$ws.Cells.Item(14, 5).Value = 0.1523527489141934
$ws.Cells.Item(14, 6).Value = 1.278846356141145
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 10).Value = 0.1577389708635693
$ws.Cells.Item(14, 13).Value = 0.4943289943925038
$ws.Cells.Item(14, 15).Value = 2.973777384702942

# Row 15
$ws.Cells.Item(15, 2).Value = 1.431526000947258
$ws.Cells.Item(15, 3).Value = 0.2972742761231189
$ws.Cells.Item(15, 4).Value = 0.1984655135534581
$ws.Cells.Item(15, 5).Value = 0.1522903504687534
$ws.Cells.Item(15, 6).Value = 1.278925422222954
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 10).Value = 0.1578255429373243
$ws.Cells.Item(15, 13).Value = 0.4907270717080081
$ws.Cells.Item(15, 15).Value = 2.974928368624546

# Row 16
$ws.Cells.Item(16, 2).Value = 1.35282670221261
$ws.Cells.Item(16, 3).Value = 0.2790463903797047
$ws.Cells.Item(16, 4).Value = 0.196467485356024
$ws.Cells.Item(16, 5).Value = 0.1519524724046875
$ws.Cells.Item(16, 6).Value = 1.27965486864862
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 10).Value = 0.1583487707997655
$ws.Cells.Item(16, 13).Value = 0.4701029273357378
$ws.Cells.Item(16, 15).Value = 2.982251925700723

# Row 17
$ws.Cells.Item(17, 2).Value = 1.304533411396392
$ws.Cells.Item(17, 3).Value = 0.2678470870287128
$ws.Cells.Item(17, 4).Value = 0.1952555336560096
$ws.Cells.Item(17, 5).Value = 0.1517627889971571
$ws.Cells.Item(17, 6).Value = 1.280348694960324
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 10).Value = 0.1586939392746523
$ws.Cells.Item(17, 13).Value = 0.4574661695717523
$ws.Cells.Item(17, 15).Value = 2.987393061214647

# Row 18
$ws.Cells.Item(18, 2).Value = 1.276751764928576
$ws.Cells.Item(18, 3).Value = 0.261399386508657
$ws.Cells.Item(18, 4).Value = 0.1945635346280739
$ws.Cells.Item(18, 5).Value = 0.1516601590391033
$ws.Cells.Item(18, 6).Value = 1.280838190060891
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 10).Value = 0.158901353837166
$ws.Cells.Item(18, 13).Value = 0.450203651861564
$ws.Cells.Item(18, 15).Value = 2.99058798352317

# Row 19
$ws.Cells.Item(19, 2).Value = 1.267344626225338
$ws.Cells.Item(19, 3).Value = 0.2592152599331143
$ws.Cells.Item(19, 4).Value = 0.1943301108111228
$ws.Cells.Item(19, 5).Value = 0.1516265221835802
$ws.Cells.Item(19, 6).Value = 1.281019445458895
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 10).Value = 0.1589731059811221
$ws.Cells.Item(19, 13).Value = 0.4477456976213716
$ws.Cells.Item(19, 15).Value = 2.991710547283958

# Row 20
$ws.Cells.Item(20, 2).Value = 1.309674803413486
$ws.Cells.Item(20, 3).Value = 0.2690399120120617
$ws.Cells.Item(20, 4).Value = 0.1953840223468717
$ws.Cells.Item(20, 5).Value = 0.1517823115220978
$ws.Cells.Item(20, 6).Value = 1.280265475044303
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 10).Value = 0.1586562760495625
$ws.Cells.Item(20, 13).Value = 0.4588107752877093
$ws.Cells.Item(20, 15).Value = 2.986821150206424

# Row 21
$ws.Cells.Item(21, 2).Value = 1.451839970193419
$ws.Cells.Item(21, 3).Value = 0.3019749770461715
$ws.Cells.Item(21, 4).Value = 0.1989856058711439
$ws.Cells.Item(21, 5).Value = 0.1523830189958879
$ws.Cells.Item(21, 6).Value = 1.278813333869294
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 10).Value = 0.1576979370879243
$ws.Cells.Item(21, 13).Value = 0.4960565110206687
$ws.Cells.Item(21, 15).Value = 2.97323834374032

# Row 22
$ws.Cells.Item(22, 2).Value = 1.54470230992888
$ws.Cells.Item(22, 3).Value = 0.3234425436016579
$ws.Cells.Item(22, 4).Value = 0.2013842577541283
$ws.Cells.Item(22, 5).Value = 0.1528331276823423
$ws.Cells.Item(22, 6).Value = 1.278671672062671
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 10).Value = 0.1571508172308853
$ws.Cells.Item(22, 13).Value = 0.5204479641238322
$ws.Cells.Item(22, 15).Value = 2.966487884969268

# Row 23
$ws.Cells.Item(23, 2).Value = 1.495145045117056
$ws.Cells.Item(23, 3).Value = 0.3119902751754182
$ws.Cells.Item(23, 4).Value = 0.2000999742361387
$ws.Cells.Item(23, 5).Value = 0.152587639250175
$ws.Cells.Item(23, 6).Value = 1.278673297663275
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 10).Value = 0.1574355723417824
$ws.Cells.Item(23, 13).Value = 0.5074254056939651
$ws.Cells.Item(23, 15).Value = 2.969895732445877

# Row 24
$ws.Cells.Item(24, 2).Value = 1.307350433524277
$ws.Cells.Item(24, 3).Value = 0.2685006640692791
$ws.Cells.Item(24, 4).Value = 0.1953259177529532
$ws.Cells.Item(24, 5).Value = 0.1517734653850162
$ws.Cells.Item(24, 6).Value = 1.280302816559256
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 10).Value = 0.1586732756458744
$ws.Cells.Item(24, 13).Value = 0.458202871159358
$ws.Cells.Item(24, 15).Value = 2.987078966113643

# Row 25
$ws.Cells.Item(25, 2).Value = 1.104611925491554
$ws.Cells.Item(25, 3).Value = 0.221356066881981
$ws.Cells.Item(25, 4).Value = 0.1903709329043863
$ws.Cells.Item(25, 5).Value = 0.1511427021034528
$ws.Cells.Item(25, 6).Value = 1.285516265775058
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 10).Value = 0.1603478792324395
$ws.Cells.Item(25, 13).Value = 0.4053319083973648
$ws.Cells.Item(25, 15).Value = 3.014703246928178
